# PG2016-1-CDP.xlsx update
# - CDP sheet (Tabela3, activities table):
#     * AT24 (row 24) Status -> "Concluída"
#     * AT26 (row 26) Status -> "Iniciada"
#     * AT28 (row 28) Inicio/Término shifted one day later
#     * Totals row: "Término" totals label "29 dias" -> "30 dias"
#     * Totals row: Status formula now explicitly counts "<>Concluída" instead of blanks
#     * Selection moved to E30

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDP")

# AT24 - mark as Concluída
$ws.Range("J24").Value2 = "Concluída"

# AT26 - mark as Iniciada
$ws.Range("J26").Value2 = "Iniciada"

# AT28 - Inicio/Término move forward one day (04/11 - 05/11 -> 05/11 - 06/11)
$ws.Range("D28").Value2 = 42679
$ws.Range("E28").Value2 = 42680

# Totals row ("29 dias" -> "30 dias")
$ws.Range("E29").Value2 = "30 dias"

# Totals row Status formula: exclude non-"Concluída" explicitly rather than blanks
$ws.Range("J29").Formula = '=(COUNTIF(J3:J28,"Concluída")/(COUNTIF(J3:J28,"Concluída") + COUNTIF(J3:J28,"<>Concluída")))'

# Move active selection to E30
$ws.Range("E30").Select()
